$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (2023-10-03 -> 2023-10-04, i.e. 45202 -> 45203) for every data row (2..135).
for ($r = 2; $r -le 135; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}
